$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.894.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.951.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +10.38%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +5.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.944.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.451"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "  +5.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.01%  "
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.438.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.947.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.50%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.838.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "416.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.699"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.11%  "
$ws.Range("E24").Value = "  +4.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.40%  "
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("E29").Value = "  +6.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0966"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.97%  "
$ws.Range("E34").Value = "  +6.61%  "
$ws.Range("E35").Value = "  +6.12%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0702"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +14.89%  "
$ws.Range("B38").Value = "Cosmos"
$ws.Range("C38").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.15"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +17.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "382.64"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.93%  "
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.709.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.238"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.52%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("E51").Value = "  +3.91%  "
